$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186296701431274
$ws.Range("B1").Value = 2.424749612808228
$ws.Range("C1").Value = 3.893263101577759
$ws.Range("D1").Value = 2.136294603347778
$ws.Range("E1").Value = 1.196138024330139
